$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '42.981.96'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = "'" + '2.542.28'
$ws.Range("E3").Value = '  +0.93%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = "'" + '317.12'
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = "'" + '97.66'
$ws.Range("E6").Value = '  +2.93%  '
$ws.Range("E7").Value = '  -0.65%  '
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").Value = "'" + '36.11'
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  +0.45%  '
$ws.Range("D12").Value = "'" + '7.64'
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("D14").Value = "'" + '2.927.13'
$ws.Range("E14").Value = '  +0.75%  '
$ws.Range("D15").Value = "'" + '2.503.06'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = "'" + '15.18'
$ws.Range("E16").Value = '  -2.55%  '
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = "'" + '42.999.74'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D20").Value = "'" + '12.84'
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("D21").Value = "'" + '0.0₃0966'
$ws.Range("E21").Value = '  -0.27%  '
$ws.Range("D22").Value = "'" + '69.74'
$ws.Range("D23").Value = "'" + '253.62'
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").Value = "'" + '2.96'
$ws.Range("E24").Value = '  -0.40%  '
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = "'" + '26.51'
$ws.Range("E26").Value = '  -3.30%  '
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("D28").Value = "'" + '2.41'
$ws.Range("E28").Value = '  +3.39%  '
$ws.Range("D29").Value = "'" + '40.89'
$ws.Range("E29").Value = '  +4.39%  '
$ws.Range("D30").Value = "'" + '10.45'
$ws.Range("E30").Value = '  +3.79%  '
$ws.Range("E31").Value = '  +0.37%  '
$ws.Range("D32").Value = "'" + '157.84'
$ws.Range("E32").Value = '  +0.84%  '
$ws.Range("D33").Value = "'" + '2.18'
$ws.Range("E33").Value = '  +4.10%  '
$ws.Range("E34").Value = '  +1.09%  '
$ws.Range("D35").Value = "'" + '2.72'
$ws.Range("E35").Value = '  +4.68%  '
$ws.Range("D36").Value = "'" + '19.10'
$ws.Range("E36").Value = '  -4.15%  '
$ws.Range("D37").Value = "'" + '0.0792'
$ws.Range("E37").Value = '  +0.94%  '
$ws.Range("E38").Value = '  +0.18%  '
$ws.Range("D39").Value = "'" + '2.49'
$ws.Range("E39").Value = '  +15.61%  '
$ws.Range("D40").Value = "'" + '0.119'
$ws.Range("E40").Value = '  -0.78%  '
$ws.Range("D41").Value = "'" + '21.89'
$ws.Range("E41").Value = '  -10.22%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("D45").Value = "'" + '3.30'
$ws.Range("E45").Value = '  -2.22%  '
$ws.Range("D46").Value = "'" + '2.021.73'
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("D47").Value = "'" + '9.09'
$ws.Range("E47").Value = '  +3.16%  '
$ws.Range("D48").Value = "'" + '84.55'
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("D49").Value = "'" + '77.84'
$ws.Range("E49").Value = '  +4.56%  '
$ws.Range("D50").Value = "'" + '106.55'
$ws.Range("E50").Value = '  +4.87%  '
$ws.Range("D51").Value = "'" + '2.779.17'
$ws.Range("E51").Value = '  +0.58%  '
